$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update regression naming from "AbrilDos" to "JunioUno" for test account data
$ws.Range("F2").Value = "AnswRegrJunioUno"
$ws.Range("G2").Value = "AnsRegrJunioUno"
$ws.Range("F3").Value = "MattioliRegrJunioUno"
$ws.Range("G3").Value = "PruebaRegrJunioUno"

# Update account creation dates and sequential IDs for pre-prod regression R33
$ws.Range("H2").Value = 20300128
$ws.Range("O2").Value = 126

$ws.Range("H3").Value = 20300128
$ws.Range("O3").Value = 127

# Update the active selection to O4
$ws.Activate()
$ws.Range("O4").Select()
